$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.229551
$ws.Range("N2").Value = 0.6886530000000001
$ws.Range("O2").Value = 0.0620874463138416
$ws.Range("P2").Value = 0.06432727452414699
$ws.Range("Q2").Value = 0.03481156218400001
$ws.Range("R2").Value = 0.3133040596560001
$ws.Range("S2").Value = 0.0620874463138416
$ws.Range("T2").Value = 0.06432727452414699

# Row 3
$ws.Range("O3").Value = 0.8197852215571831
$ws.Range("P3").Value = 0.8493592848284227
$ws.Range("S3").Value = 0.8197852215571831
$ws.Range("T3").Value = 0.8493592848284227

# Row 4
$ws.Range("M4").Value = 0.01067566666666667
$ws.Range("N4").Value = 0.032027
$ws.Range("O4").Value = 0.00288748418012178
$ws.Range("P4").Value = 0.002991651268759238
$ws.Range("Q4").Value = 0.001618971967111111
$ws.Range("R4").Value = 0.014570747704
$ws.Range("S4").Value = 0.00288748418012178
$ws.Range("T4").Value = 0.002991651268759238

# Row 5
$ws.Range("M5").Value = 0.3862035
$ws.Range("N5").Value = 0.772407
$ws.Range("O5").Value = 0.1044577852959374
$ws.Range("P5").Value = 0.0721507597198775
$ws.Range("Q5").Value = 0.058568018244
$ws.Range("R5").Value = 0.351408109464
$ws.Range("S5").Value = 0.1044577852959374
$ws.Range("T5").Value = 0.0721507597198775

# Row 6
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03986366666666667
$ws.Range("N6").Value = 0.119591
$ws.Range("O6").Value = 0.0107820626529161
$ws.Range("P6").Value = 0.01117102965879371
$ws.Range("Q6").Value = 0.006045351625777779
$ws.Range("R6").Value = 0.054408164632
$ws.Range("S6").Value = 0.0107820626529161
$ws.Range("T6").Value = 0.01117102965879371
